$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 261393.42
$ws.Range("I113").Value = 385872.78
$ws.Range("K113").Value = 385872.78
$ws.Range("M113").Value = -382618.78
$ws.Range("H132").Value = 1303.4688
$ws.Range("I132").Value = 815.1786
$ws.Range("J132").Value = 4721.5
$ws.Range("K132").Value = 2445.5358
$ws.Range("L132").Value = 14164.5
$ws.Range("M132").Value = 84.46420000000035
$ws.Range("N132").Value = -19224.5
$ws.Range("H137").Value = 1086.9524
$ws.Range("I137").Value = 1021.86664
$ws.Range("K137").Value = 3065.59992
$ws.Range("M137").Value = -515.5999199999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 6780
$ws.Range("I21").Value = 3633.3333
$ws.Range("J21").Value = 11500
$ws.Range("K21").Value = 3633.3333
$ws.Range("L21").Value = 11500
$ws.Range("M21").Value = -3259.3333
$ws.Range("N21").Value = -12248

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 46000
$ws.Range("J2").Value = 46000
$ws.Range("L2").Value = 46000
$ws.Range("N2").Value = -46226
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H134").Value = 2253.8408
$ws.Range("I134").Value = 2044.7742
$ws.Range("K134").Value = 6134.3226
$ws.Range("M134").Value = -3599.3226

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 79.583336
$ws.Range("I7").Value = 86.5
$ws.Range("K7").Value = 86.5
$ws.Range("M7").Value = 26.5
$ws.Range("H55").Value = 14593
$ws.Range("I55").Value = 9900
$ws.Range("J55").Value = 23979
$ws.Range("K55").Value = 9900
$ws.Range("L55").Value = 23979
$ws.Range("M55").Value = -9585
$ws.Range("N55").Value = -24609
$ws.Range("H69").Value = 24100.5
$ws.Range("J69").Value = 24100.5
$ws.Range("L69").Value = 24100.5
$ws.Range("N69").Value = -25598.5
$ws.Range("H72").Value = 24100.5
$ws.Range("J72").Value = 24100.5
$ws.Range("L72").Value = 72301.5
$ws.Range("N72").Value = -79789.5
$ws.Range("H105").Value = 1490
$ws.Range("I105").Value = 1538.3334
$ws.Range("J105").Value = 1200
$ws.Range("K105").Value = 1538.3334
$ws.Range("L105").Value = 1200
$ws.Range("M105").Value = 208.6666
$ws.Range("N105").Value = -4694
$ws.Range("H132").Value = 3558.8
$ws.Range("I132").Value = 2600
$ws.Range("J132").Value = 4997
$ws.Range("K132").Value = 7800
$ws.Range("L132").Value = 14991
$ws.Range("M132").Value = -5270
$ws.Range("N132").Value = -20051

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1275.5714
$ws.Range("I4").Value = 227.25
$ws.Range("J4").Value = 4630.2
$ws.Range("K4").Value = 681.75
$ws.Range("L4").Value = 13890.6
$ws.Range("M4").Value = -569.75
$ws.Range("N4").Value = -14114.6
$ws.Range("H6").Value = 89.80952499999999
$ws.Range("I6").Value = 64.3
$ws.Range("K6").Value = 192.9
$ws.Range("M6").Value = -79.89999999999998
$ws.Range("H10").Value = 37
$ws.Range("I10").Value = 37
$ws.Range("K10").Value = 111
$ws.Range("M10").Value = 28
$ws.Range("H122").Value = 4650.92
$ws.Range("I122").Value = 547.1111
$ws.Range("J122").Value = 6959.3125
$ws.Range("K122").Value = 4923.9999
$ws.Range("L122").Value = 62633.8125
$ws.Range("M122").Value = -2473.9999
$ws.Range("N122").Value = -67533.8125
$ws.Range("H125").Value = 9371.5
$ws.Range("I125").Value = 1500
$ws.Range("J125").Value = 11995.333
$ws.Range("K125").Value = 4500
$ws.Range("L125").Value = 35985.999
$ws.Range("M125").Value = 420
$ws.Range("N125").Value = -45825.999
$ws.Range("H134").Value = 3907.4905
$ws.Range("I134").Value = 1373.1333
$ws.Range("J134").Value = 4907.8945
$ws.Range("K134").Value = 4119.3999
$ws.Range("L134").Value = 14723.6835
$ws.Range("M134").Value = 950.6000999999997
$ws.Range("N134").Value = -24863.6835
$ws.Range("H141").Value = 4145.0713
$ws.Range("I141").Value = 3374.8333
$ws.Range("K141").Value = 10124.4999
$ws.Range("M141").Value = -4944.499899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 23320
$ws.Range("J51").Value = 23320
$ws.Range("L51").Value = 23320
$ws.Range("N51").Value = -24338

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2122.4443
$ws.Range("I7").Value = 1912.75
$ws.Range("J7").Value = 3800
$ws.Range("K7").Value = 1912.75
$ws.Range("L7").Value = 3800
$ws.Range("M7").Value = -1800.75
$ws.Range("N7").Value = -4024
$ws.Range("H46").Value = 1371.5714
$ws.Range("I46").Value = 1340.2
$ws.Range("J46").Value = 1450
$ws.Range("K46").Value = 1340.2
$ws.Range("L46").Value = 1450
$ws.Range("M46").Value = -1152.2
$ws.Range("N46").Value = -1826
$ws.Range("H55").Value = 280.06668
$ws.Range("I55").Value = 209
$ws.Range("J55").Value = 386.66666
$ws.Range("K55").Value = 209
$ws.Range("L55").Value = 386.66666
$ws.Range("M55").Value = -36
$ws.Range("N55").Value = -732.66666
$ws.Range("H68").Value = 8126.6665
$ws.Range("I68").Value = 21120
$ws.Range("J68").Value = 1630
$ws.Range("K68").Value = 21120
$ws.Range("L68").Value = 1630
$ws.Range("M68").Value = -20371
$ws.Range("N68").Value = -3128
$ws.Range("H71").Value = 8126.6665
$ws.Range("I71").Value = 21120
$ws.Range("J71").Value = 1630
$ws.Range("K71").Value = 105600
$ws.Range("L71").Value = 8150
$ws.Range("M71").Value = -101856
$ws.Range("N71").Value = -15638
$ws.Range("H93").Value = 5916.6816
$ws.Range("I93").Value = 10798.6
$ws.Range("J93").Value = 1848.4166
$ws.Range("K93").Value = 10798.6
$ws.Range("L93").Value = 1848.4166
$ws.Range("M93").Value = -9550.6
$ws.Range("N93").Value = -4344.4166
$ws.Range("H122").Value = 14867357
$ws.Range("I122").Value = 1854255.5
$ws.Range("J122").Value = 50002730
$ws.Range("K122").Value = 5562766.5
$ws.Range("L122").Value = 150008190
$ws.Range("M122").Value = -5560316.5
$ws.Range("N122").Value = -150013090
$ws.Range("H126").Value = 2122.4443
$ws.Range("I126").Value = 1912.75
$ws.Range("J126").Value = 3800
$ws.Range("K126").Value = 5738.25
$ws.Range("L126").Value = 11400
$ws.Range("M126").Value = -3268.25
$ws.Range("N126").Value = -16340
$ws.Range("H132").Value = 2648.52
$ws.Range("I132").Value = 1883.2354
$ws.Range("J132").Value = 4274.75
$ws.Range("K132").Value = 5649.706200000001
$ws.Range("L132").Value = 12824.25
$ws.Range("M132").Value = -3119.706200000001
$ws.Range("N132").Value = -17884.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 4900
$ws.Range("J26").Value = 4900
$ws.Range("L26").Value = 4900
$ws.Range("N26").Value = -5486
$ws.Range("H96").Value = 3400
$ws.Range("J96").Value = 3400
$ws.Range("L96").Value = 3400
$ws.Range("N96").Value = -6146
$ws.Range("H100").Value = 5051274.5
$ws.Range("I100").Value = 5682559
$ws.Range("K100").Value = 11365118
$ws.Range("M100").Value = -11364577
